$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=45689.55416666667; B=99.44},
    @{Row=3;  A=45692.07638888889; B=87.98},
    @{Row=4;  A=45692.57916666667; B=100},
    @{Row=5;  A=45693.60833333333; B=99.28},
    @{Row=6;  A=45695.93611111111; B=95.12},
    @{Row=7;  A=45699.20416666667; B=62.42},
    @{Row=8;  A=45703.20277777778; B=62.49},
    @{Row=9;  A=45707.81736111111; B=100},
    @{Row=10; A=45711.36875;       B=69.79000000000001},
    @{Row=11; A=45721.63125;       B=99.41}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}
